$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header cells: "_old" -> "_FV2404", "_new" -> "_FV2410"
$headers = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Stash the header row's existing formatting (bold/fill/borders/wrap/center)
#    in a scratch range so we can restore it verbatim after the table is
#    created - creating a ListObject auto-applies its own header dxf/style,
#    which we do not want, so the header range is cleared before Add() and
#    the original look is pasted back afterwards.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null
$headerRange.ClearFormats() | Out-Null

# 3. Turn the A1:U78 range into a table ("Table1") with the (new) header names
$tableRange = $ws.Range("A1:U78")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Restore the header formatting captured in step 2, then clean up the scratch area
$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null
$scratch.ClearFormats() | Out-Null
$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = $false

# 4. Freeze the header row (split after row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null

Write-Output "Header row renamed, Table1 added, header row frozen."
